$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = -0.06184127475628475
$ws.Range("C2").Value = 0.6728703294834245
$ws.Range("D2").Value = 1.036433932025511
$ws.Range("E2").Value = 1.018053992686788
$ws.Range("F2").Value = 1.040087243427768
$ws.Range("G2").Value = 22

# Row 3 (Q0)
$ws.Range("B3").Value = 0.1029130003429865
$ws.Range("C3").Value = 0.8770720106948019
$ws.Range("D3").Value = 1.404277130504535
$ws.Range("E3").Value = 1.185021995789334
$ws.Range("F3").Value = 1.207077218135037
$ws.Range("G3").Value = 23

# Row 4 (Q1)
$ws.Range("B4").Value = 0.1140919098201663
$ws.Range("C4").Value = 1.461248407423022
$ws.Range("D4").Value = 9.665212501827879
$ws.Range("E4").Value = 3.108892487981513
$ws.Range("F4").Value = 3.179909416460013
$ws.Range("G4").Value = 22
